# Branched from CRAN release
# Adds a new "2.5.0" benchmark row (row 27) to the httk-benchmarks table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The data lives in an Excel Table (ListObject) named "Table1" which spans
# A1:R26. Add a new table row - this automatically grows the table range
# (and the AutoFilter range) from A1:R26 to A1:R27.
$tbl = $ws.ListObjects.Item(1)
$newListRow = $tbl.ListRows.Add()

# Match the left-aligned formatting used by every other data row.
$ws.Range("A27:R27").HorizontalAlignment = -4131

# Populate the new row - essentially a re-run of the previous release's
# (2.4.0, row 26) benchmark numbers under the new 2.5.0 version tag, with a
# small change to the E column and a new release note.
$ws.Range("A27").Value = "2.5.0"
$ws.Range("B27").Value = 1021
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0.9999
$ws.Range("F27").Value = 0.9477
$ws.Range("G27").Value = 353
$ws.Range("H27").Value = 0.2716
$ws.Range("I27").Value = 353
$ws.Range("J27").Value = 1.508
$ws.Range("K27").Value = 36
$ws.Range("L27").Value = 0.9698
$ws.Range("M27").Value = 80
$ws.Range("N27").Value = 1.132
$ws.Range("O27").Value = 80
$ws.Range("P27").Value = 0.6466
$ws.Range("Q27").Value = 863
$ws.Range("R27").Value = "Added models 3comp2 and sumclearances"

# Move the view/selection to the newly-added last cell, mirroring where a
# person would naturally end up after typing in the new row.
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 16
[void]$ws.Range("R27").Select()

Write-Host "Added 2.5.0 benchmark row to Sheet1/Table1"
